# Scheduled market-data refresh for the Halicarnassus Profits workbook.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) per leve row
# with freshly pulled market values; a few rows also gain/lose an HQ-profit (N) or
# NQ-profit (M) cell depending on whether that recipe currently has an HQ/NQ variant.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 827.3333
$ws.Range("I4").Value = 741.25
$ws.Range("K4").Value = 741.25
$ws.Range("M4").Value = -627.25
# Row 5
$ws.Range("H5").Value = 134.8
$ws.Range("J5").Value = 143.75
$ws.Range("L5").Value = 143.75
$ws.Range("N5").Value = -373.75
# Row 9
$ws.Range("H9").Value = 190.36363
$ws.Range("I9").Value = 204.14285
$ws.Range("K9").Value = 204.14285
$ws.Range("M9").Value = -35.14285000000001
# Row 55
$ws.Range("H55").Value = 348.3158
$ws.Range("I55").Value = 367.77777
$ws.Range("J55").Value = 330.8
$ws.Range("K55").Value = 367.77777
$ws.Range("L55").Value = 330.8
$ws.Range("M55").Value = -153.77777
$ws.Range("N55").Value = -758.8
# Row 82
$ws.Range("H82").Value = 699
$ws.Range("I82").Value = 699
$ws.Range("K82").Value = 2097
$ws.Range("M82").Value = -1691
# Row 85
$ws.Range("H85").Value = 699
$ws.Range("I85").Value = 699
$ws.Range("K85").Value = 2097
$ws.Range("M85").Value = -693
# Row 95
$ws.Range("H95").Value = 26540.666
$ws.Range("J95").Value = 26540.666
$ws.Range("L95").Value = 26540.666
$ws.Range("N95").Value = -32032.666
# Row 101
$ws.Range("I101").Value = 248
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 744
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 878
$ws.Range("N101").ClearContents()
# Row 116
$ws.Range("H116").Value = 3946.6667
$ws.Range("I116").Value = 3946.6667
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 3946.6667
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -504.6667000000002
$ws.Range("N116").ClearContents()
# Row 127
$ws.Range("H127").Value = 624
$ws.Range("I127").Value = 624
$ws.Range("K127").Value = 1872
$ws.Range("M127").Value = 3088
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 3279.5
$ws.Range("I22").Value = 3366.111
$ws.Range("K22").Value = 3366.111
$ws.Range("M22").Value = -3193.111

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1417.8462
$ws.Range("I22").Value = 1238.3334
$ws.Range("J22").Value = 1571.7142
$ws.Range("K22").Value = 1238.3334
$ws.Range("L22").Value = 1571.7142
$ws.Range("M22").Value = -888.3334
$ws.Range("N22").Value = -2271.7142
# Row 58
$ws.Range("H58").Value = 5031.909
$ws.Range("I58").Value = 4550.75
$ws.Range("K58").Value = 4550.75
$ws.Range("M58").Value = -4347.75
# Row 136
$ws.Range("H136").Value = 5031.909
$ws.Range("I136").Value = 4550.75
$ws.Range("K136").Value = 13652.25
$ws.Range("M136").Value = -11102.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 484.75
$ws.Range("I14").Value = 484.75
$ws.Range("K14").Value = 1454.25
$ws.Range("M14").Value = -1281.25
# Row 40
$ws.Range("H40").Value = 118.8125
$ws.Range("I40").Value = 25.083334
$ws.Range("K40").Value = 100.333336
$ws.Range("M40").Value = -31.333336
# Row 49
$ws.Range("H49").Value = 1000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 1000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 3000
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -3312
# Row 129
$ws.Range("H129").Value = 1000
$ws.Range("I129").Value = 1000
$ws.Range("K129").Value = 3000
$ws.Range("M129").Value = 2000

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 297.91306
$ws.Range("I2").Value = 171.14285
$ws.Range("J2").Value = 495.1111
$ws.Range("K2").Value = 171.14285
$ws.Range("L2").Value = 495.1111
$ws.Range("M2").Value = -58.14285000000001
$ws.Range("N2").Value = -721.1111000000001
# Row 80
$ws.Range("H80").Value = 2364.8333
$ws.Range("I80").Value = 2794.5
$ws.Range("J80").Value = 2150
$ws.Range("K80").Value = 2794.5
$ws.Range("L80").Value = 2150
$ws.Range("M80").Value = -1796.5
$ws.Range("N80").Value = -4146
# Row 83
$ws.Range("H83").Value = 2364.8333
$ws.Range("I83").Value = 2794.5
$ws.Range("J83").Value = 2150
$ws.Range("K83").Value = 13972.5
$ws.Range("L83").Value = 10750
$ws.Range("M83").Value = -8980.5
$ws.Range("N83").Value = -20734
# Row 102
$ws.Range("H102").Value = 2353.7273
$ws.Range("I102").Value = 1876.7778
$ws.Range("K102").Value = 1876.7778
$ws.Range("M102").Value = -254.7778000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7811.615
$ws.Range("J7").Value = 8198.799999999999
$ws.Range("L7").Value = 8198.799999999999
$ws.Range("N7").Value = -8422.799999999999
# Row 19
$ws.Range("H19").Value = 3500
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 3500
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 3500
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -3840
# Row 46
$ws.Range("H46").Value = 6617.0586
$ws.Range("J46").Value = 6820.7144
$ws.Range("L46").Value = 6820.7144
$ws.Range("N46").Value = -7196.7144
# Row 55
$ws.Range("H55").Value = 1013.3571
$ws.Range("I55").Value = 1084
$ws.Range("J55").Value = 836.75
$ws.Range("K55").Value = 1084
$ws.Range("L55").Value = 836.75
$ws.Range("M55").Value = -911
$ws.Range("N55").Value = -1182.75
# Row 68
$ws.Range("H68").Value = 7393.2856
$ws.Range("I68").Value = 4084.3333
$ws.Range("K68").Value = 4084.3333
$ws.Range("M68").Value = -3335.3333
# Row 71
$ws.Range("H71").Value = 7393.2856
$ws.Range("I71").Value = 4084.3333
$ws.Range("K71").Value = 20421.6665
$ws.Range("M71").Value = -16677.6665
# Row 76
$ws.Range("H76").Value = 19993.5
$ws.Range("J76").Value = 19993.5
$ws.Range("L76").Value = 19993.5
$ws.Range("N76").Value = -20669.5
# Row 79
$ws.Range("H79").Value = 19993.5
$ws.Range("J79").Value = 19993.5
$ws.Range("L79").Value = 19993.5
$ws.Range("N79").Value = -22333.5
# Row 82
$ws.Range("H82").Value = 2775.3
$ws.Range("I82").Value = 1419.2
$ws.Range("K82").Value = 1419.2
$ws.Range("M82").Value = -1058.2
# Row 85
$ws.Range("H85").Value = 2775.3
$ws.Range("I85").Value = 1419.2
$ws.Range("K85").Value = 1419.2
$ws.Range("M85").Value = -171.2
# Row 126
$ws.Range("H126").Value = 7811.615
$ws.Range("J126").Value = 8198.799999999999
$ws.Range("L126").Value = 24596.4
$ws.Range("N126").Value = -29536.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 10000
$ws.Range("J15").Value = 10000
$ws.Range("L15").Value = 10000
$ws.Range("N15").Value = -10576
# Row 113
$ws.Range("H113").Value = 850.9
$ws.Range("I113").Value = 751.25
$ws.Range("K113").Value = 2253.75
$ws.Range("M113").Value = -83.75
